# Refresh the cryptocurrency snapshot table (columns B-E, rows 2-51) on
# Sheet1 with the latest scraped values. Column D ('Price') holds plain
# text that merely looks numeric (e.g. "212.77", "1.639.74"); Excel would
# otherwise silently coerce such a string into a real number the moment
# .Value is assigned (dropping trailing zeros, choking on the
# thousands-style extra dot, etc.), so each Price cell is pinned to the
# Text number format first to keep it a literal string, matching the
# original inline-string cell content exactly.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '27.967.03'
$ws.Range('E2').Value = '  +0.87%  '

# Row 3
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.639.74'
$ws.Range('E3').Value = '  +0.36%  '

# Row 4
$ws.Range('E4').Value = '  +0.09%  '

# Row 5
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '212.77'
$ws.Range('E5').Value = '  +0.33%  '

# Row 6
$ws.Range('E6').Value = '  +0.39%  '

# Row 7
$ws.Range('E7').Value = '  +0.05%  '

# Row 8
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '23.46'
$ws.Range('E8').Value = '  +1.11%  '

# Row 9
$ws.Range('E9').Value = '  -2.27%  '

# Row 10
$ws.Range('E10').Value = '  +0.44%  '

# Row 11
$ws.Range('E11').Value = '  +2.30%  '

# Row 12
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '1.873.41'
$ws.Range('E12').Value = '  +0.43%  '

# Row 13
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '1.629.95'
$ws.Range('E13').Value = '  -0.32%  '

# Row 14
$ws.Range('E14').Value = '  +3.76%  '

# Row 15
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '4.09'
$ws.Range('E15').Value = '  +1.35%  '

# Row 16
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '65.78'
$ws.Range('E16').Value = '  +0.95%  '

# Row 17
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '27.966.14'
$ws.Range('E17').Value = '  +1.05%  '

# Row 18
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '234.65'
$ws.Range('E18').Value = '  +2.33%  '

# Row 19
$ws.Range('E19').Value = '  +0.64%  '

# Row 20
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '7.60'
$ws.Range('E20').Value = '  +0.47%  '

# Row 21
$ws.Range('E21').Value = '  +0.05%  '

# Row 22
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '10.59'
$ws.Range('E22').Value = '  -0.67%  '

# Row 23
$ws.Range('E23').Value = '  +0.52%  '

# Row 24
$ws.Range('E24').Value = '  -1.78%  '

# Row 25
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '151.82'

# Row 26
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '6.95'
$ws.Range('E26').Value = '  +1.29%  '

# Row 27 (Stellar and EthereumClassic swapped rank with row 28)
$ws.Range('B27').Value = 'Stellar'
$ws.Range('C27').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '0.111'
$ws.Range('E27').Value = '  +0.14%  '

# Row 28
$ws.Range('B28').Value = 'EthereumClassic'
$ws.Range('C28').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '15.66'
$ws.Range('E28').Value = '  +0.51%  '

# Row 29
$ws.Range('E29').Value = '  +0.07%  '

# Row 30
$ws.Range('E30').Value = '  +0.56%  '

# Row 31
$ws.Range('E31').Value = '  +0.54%  '

# Row 32
$ws.Range('E32').Value = '  +2.02%  '

# Row 33
$ws.Range('E33').Value = '  +1.23%  '

# Row 34
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.414.11'
$ws.Range('E34').Value = '  -3.86%  '

# Row 35
$ws.Range('E35').Value = '  +2.06%  '

# Row 36
$ws.Range('E36').Value = '  +1.42%  '

# Row 37
$ws.Range('E37').Value = '  +1.45%  '

# Row 38
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.881'
$ws.Range('E38').Value = '  +0.53%  '

# Row 39
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.557'
$ws.Range('E39').Value = '  -0.21%  '

# Row 40
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.901'
$ws.Range('E40').Value = '  -3.52%  '

# Row 41
$ws.Range('E41').Value = '  +1.05%  '

# Row 42
$ws.Range('E42').Value = '  +0.02%  '

# Row 43
$ws.Range('E43').Value = '  +6.67%  '

# Row 44
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '66.58'
$ws.Range('E44').Value = '  -1.74%  '

# Row 45
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '5.52'
$ws.Range('E45').Value = '  +3.01%  '

# Row 46
$ws.Range('E46').Value = '  -0.19%  '

# Row 47
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '1.782.45'
$ws.Range('E47').Value = '  +0.55%  '

# Row 48
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '87.72'
$ws.Range('E48').Value = '  +0.12%  '

# Row 50
$ws.Range('E50').Value = '  +0.25%  '

# Row 51
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '7.63'
$ws.Range('E51').Value = '  -1.27%  '
